# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 10991
$ws.Range("C4").Value = 3849
$ws.Range("C5").Value = 196
$ws.Range("D5").Value = 98.8
$ws.Range("C6").Value = 15036
